$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Style = "Normal"

$metaLabel = "Meta description"
$metaRest  = ": Read our review of Boots of Luck, an Irish-themed slot game by Betixon with 5*3 reel system and 20 paylines. Play for free and enjoy simple gameplay mechanics."

$metaRange = $metaPara.Range
$metaRange.Text = $metaLabel + $metaRest

# Bold just the "Meta description" label portion.
$metaStart = $metaPara.Range.Start
$boldRange = $d.Range($metaStart, $metaStart + $metaLabel.Length)
$boldRange.Font.Bold = 1

# ---------------------------------------------------------------------------
# 2) Near the bottom: drop the duplicated bold title paragraph, and replace
#    the italic "meta description" paragraph's text with the new art prompt.
# ---------------------------------------------------------------------------
$oldBoldText = "Play Boots of Luck for Free - Game Review"
$oldItalicText = "Read our review of Boots of Luck, an Irish-themed slot game by Betixon with 5*3 reel system and 20 paylines. Play for free and enjoy simple gameplay mechanics."
$newItalicText = "Create a cartoon-style image for the game ""Boots of Luck"" that features a happy Maya warrior with glasses. The image should have a green and gold color scheme, with the warrior standing on a grassy field overlooking a lake in the background. The warrior should be wearing a green tunic with a gold belt and a gold pair of lucky boots. The warrior should be holding a mug of beer with a four-leaf clover on the side in one hand and a pot full of gold coins in the other. The warrior should be depicted with a big smile on their face and wearing a pair of glasses with a reflection of the game on the lenses."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($oldBoldText) -and $p.Range.Style.NameLocal -ne "Heading 1") {
        $p.Range.Delete()
        break
    }
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith($oldItalicText)) {
        $target = $d.Range($p.Range.Start, $p.Range.Start + $oldItalicText.Length)
        $target.Text = $newItalicText
        break
    }
}
